$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.02406616871877758
$ws.Range("J2").Value = 0.02406616871877757
$ws.Range("M2").Value = 2.765491666666666
$ws.Range("N2").Value = 8.296474999999999
$ws.Range("O2").Value = 0.4885734765724882
$ws.Range("P2").Value = 0.4885734765724882
$ws.Range("Q2").Value = 0.1239456491777778
$ws.Range("R2").Value = 1.1155108426
$ws.Range("S2").Value = 0.01175809171871323
$ws.Range("T2").Value = 0.01175809171871322

# Row 3
$ws.Range("I3").Value = 0.02406616871877758
$ws.Range("J3").Value = 0.02406616871877757
$ws.Range("N3").Value = 6.670577999999999
$ws.Range("O3").Value = 0.3928255655815217
$ws.Range("P3").Value = 0.3928255655815217
$ws.Range("Q3").Value = 0.09965547061866666
$ws.Range("R3").Value = 0.8968992355679998
$ws.Range("S3").Value = 0.009453806338334128
$ws.Range("T3").Value = 0.009453806338334125

# Row 4
$ws.Range("D4").Value = "MuSCs"
$ws.Range("I4").Value = 0.02406616871877758
$ws.Range("J4").Value = 0.02406616871877757
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6111460000000001
$ws.Range("N4").Value = 1.833438
$ws.Range("O4").Value = 0.107969851984139
$ws.Range("P4").Value = 0.107969851984139
$ws.Range("Q4").Value = 0.02739074885866667
$ws.Range("R4").Value = 0.246516739728
$ws.Range("S4").Value = 0.002598420674391732
$ws.Range("T4").Value = 0.002598420674391731

# Row 5
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("I5").Value = 0.02406616871877758
$ws.Range("J5").Value = 0.02406616871877757
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.06017566666666666
$ws.Range("N5").Value = 0.180527
$ws.Range("O5").Value = 0.01063110586185116
$ws.Range("P5").Value = 0.01063110586185116
$ws.Range("Q5").Value = 0.002696993145777777
$ws.Range("R5").Value = 0.024272938312
$ws.Range("S5").Value = 0.0002558499873384953
$ws.Range("T5").Value = 0.0002558499873384952

# Row 6
$ws.Range("I6").Value = 0.8626970447097064
$ws.Range("J6").Value = 0.8626970447097063
$ws.Range("M6").Value = 2.765491666666666
$ws.Range("N6").Value = 8.296474999999999
$ws.Range("O6").Value = 0.4885734765724882
$ws.Range("P6").Value = 0.4885734765724882
$ws.Range("Q6").Value = 4.443064722922221
$ws.Range("R6").Value = 39.98758250629999
$ws.Range("S6").Value = 0.4214908943626326
$ws.Range("T6").Value = 0.4214908943626325

# Row 7
$ws.Range("I7").Value = 0.8626970447097064
$ws.Range("J7").Value = 0.8626970447097063
$ws.Range("N7").Value = 6.670577999999999
$ws.Range("O7").Value = 0.3928255655815217
$ws.Range("P7").Value = 0.3928255655815217
$ws.Range("Q7").Value = 3.572337624509332
$ws.Range("R7").Value = 32.15103862058399
$ws.Range("S7").Value = 0.3388894545135978
$ws.Range("T7").Value = 0.3388894545135976

# Row 8
$ws.Range("D8").Value = "MuSCs"
$ws.Range("I8").Value = 0.8626970447097064
$ws.Range("J8").Value = 0.8626970447097063
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.6111460000000001
$ws.Range("N8").Value = 1.833438
$ws.Range("O8").Value = 0.107969851984139
$ws.Range("P8").Value = 0.107969851984139
$ws.Range("Q8").Value = 0.9818728676293333
$ws.Range("R8").Value = 8.836855808664
$ws.Range("S8").Value = 0.09314527222446119
$ws.Range("T8").Value = 0.09314527222446115

# Row 9
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("I9").Value = 0.8626970447097064
$ws.Range("J9").Value = 0.8626970447097063
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.06017566666666666
$ws.Range("N9").Value = 0.180527
$ws.Range("O9").Value = 0.01063110586185116
$ws.Range("P9").Value = 0.01063110586185116
$ws.Range("Q9").Value = 0.09667878770622221
$ws.Range("R9").Value = 0.8701090893559998
$ws.Range("S9").Value = 0.009171423609015031
$ws.Range("T9").Value = 0.009171423609015027

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.210882
$ws.Range("H10").Value = 0.6326459999999999
$ws.Range("I10").Value = 0.113236786571516
$ws.Range("J10").Value = 0.113236786571516
$ws.Range("M10").Value = 2.765491666666666
$ws.Range("N10").Value = 8.296474999999999
$ws.Range("O10").Value = 0.4885734765724882
$ws.Range("P10").Value = 0.4885734765724882
$ws.Range("Q10").Value = 0.5831924136499999
$ws.Range("R10").Value = 5.248731722849999
$ws.Range("S10").Value = 0.05532449049114244
$ws.Range("T10").Value = 0.05532449049114242

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.210882
$ws.Range("H11").Value = 0.6326459999999999
$ws.Range("I11").Value = 0.113236786571516
$ws.Range("J11").Value = 0.113236786571516
$ws.Range("N11").Value = 6.670577999999999
$ws.Range("O11").Value = 0.3928255655815217
$ws.Range("P11").Value = 0.3928255655815217
$ws.Range("Q11").Value = 0.4689016099319999
$ws.Range("R11").Value = 4.220114489387999
$ws.Range("S11").Value = 0.04448230472958985
$ws.Range("T11").Value = 0.04448230472958983

# Row 12
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.210882
$ws.Range("H12").Value = 0.6326459999999999
$ws.Range("I12").Value = 0.113236786571516
$ws.Range("J12").Value = 0.113236786571516
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.6111460000000001
$ws.Range("N12").Value = 1.833438
$ws.Range("O12").Value = 0.107969851984139
$ws.Range("P12").Value = 0.107969851984139
$ws.Range("Q12").Value = 0.128879690772
$ws.Range("R12").Value = 1.159917216948
$ws.Range("S12").Value = 0.01222615908528613
$ws.Range("T12").Value = 0.01222615908528612

# Row 13
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.210882
$ws.Range("H13").Value = 0.6326459999999999
$ws.Range("I13").Value = 0.113236786571516
$ws.Range("J13").Value = 0.113236786571516
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.06017566666666666
$ws.Range("N13").Value = 0.180527
$ws.Range("O13").Value = 0.01063110586185116
$ws.Range("P13").Value = 0.01063110586185116
$ws.Range("Q13").Value = 0.012689964938
$ws.Range("R13").Value = 0.114209684442
$ws.Range("S13").Value = 0.001203832265497633
$ws.Range("T13").Value = 0.001203832265497632
